# 036 Week 20/21 Update
# Fills in Week 20 scores (column X on Sheet1, column U on THURSDAY SINGLES)
# for players who were missing a result, and corrects a handicap value on
# the HANDICAPS sheet. Dependent formulas (SUM / COUNTIF) recalculate
# automatically.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: WK 20 column (X) ------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

$sheet1.Range("X12").Value = 33.0
$sheet1.Range("X13").Value = 31.0
$sheet1.Range("X20").Value = 29.0
$sheet1.Range("X21").Value = 33.0
$sheet1.Range("X22").Value = 32.0
$sheet1.Range("X26").Value = 35.0
$sheet1.Range("X27").Value = 35.0
$sheet1.Range("X28").Value = 35.0
$sheet1.Range("X31").Value = 29.0

# ---- THURSDAY SINGLES: WK 20 column (U) ---------------------------------
$sheet2 = $wb.Worksheets.Item("THURSDAY SINGLES")

$sheet2.Range("U5").Value = 37.0
$sheet2.Range("U6").Value = 32.0
$sheet2.Range("U7").Value = 27.0
$sheet2.Range("U12").Value = 32.0
$sheet2.Range("U13").Value = 34.0
$sheet2.Range("U15").Value = 23.0
$sheet2.Range("U16").Value = 27.0

# ---- HANDICAPS: correct JOHN ANTCLIFFE's handicap ------------------------
$sheet3 = $wb.Worksheets.Item("HANDICAPS")

$sheet3.Range("B7").Value = 12.0
$sheet3.Range("C7").Value = 12.0
